# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Both sheets contain the same data, so the same cell updates are applied to each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    "F2"  = 1165
    "F3"  = 89
    "F4"  = 1524
    "F7"  = 11171
    "F8"  = 85
    "F10" = 328
    "F12" = 767
    "F13" = 12259
    "F14" = 12850
    "F21" = 56
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
